# reference.docx: Include 'FootnoteText' style.
#
# Adds the built-in "Footnote Text" paragraph style (styleId FootnoteText)
# to the document's style sheet, based on Normal, with NextParagraphStyle
# pointing back at itself, uiPriority 9, unhideWhenUsed and qFormat set -
# mirroring what Word itself writes out when a FootnoteText style is first
# materialized.

$d = $word.ActiveDocument

# `Styles("Footnote Text")` / `Paragraph.Style = "Footnote Text"` both
# resolve Word's built-in style definition lazily; actually assigning it to
# a paragraph is what makes Word materialize a real (non-custom) entry with
# styleId "FootnoteText" in styles.xml. Do that on a throwaway paragraph so
# the visible document content is left untouched, then remove the
# paragraph again once the style exists.

$tail = $d.Content
$tail.Collapse(0)
$scratch = $d.Paragraphs.Add($tail)
$scratch.Style = "Footnote Text"

$style = $d.Styles("Footnote Text")
$style.BaseStyle = $d.Styles("Normal")
$style.NextParagraphStyle = "FootnoteText"
$style.Priority = 9
$style.UnhideWhenUsed = $true
$style.QuickStyle = $true

$scratch.Range.Delete()
